$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Rushing) ---
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "D.Prescott"
$ws1.Range("C2").Value = 11
$ws1.Range("D2").Value = 4
$ws1.Range("E2").Value = 10
$ws1.Range("F2").Value = 9

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "C.Rush"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "E.Elliott"
$ws1.Range("C4").Value = 112
$ws1.Range("D4").Value = 51
$ws1.Range("E4").Value = 22
$ws1.Range("F4").Value = 31

$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "T.Pollard"
$ws1.Range("C5").Value = 66
$ws1.Range("D5").Value = 33
$ws1.Range("E5").Value = 7
$ws1.Range("F5").Value = 14

$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "C.Clement"
$ws1.Range("C6").Value = 10
$ws1.Range("D6").Value = 4
$ws1.Range("E6").Value = 2
$ws1.Range("F6").Value = 2

$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "J.Hardy"
$ws1.Range("C7").Value = 1
$ws1.Range("D7").Value = 0
$ws1.Range("E7").Value = 0
$ws1.Range("F7").Value = 0

$ws1.Range("A8").Value = 6
$ws1.Range("B8").Value = "C.Lamb"
$ws1.Range("C8").Value = 3
$ws1.Range("D8").Value = 3
$ws1.Range("E8").Value = 2
$ws1.Range("F8").Value = 2

$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "C.Wilson"
$ws1.Range("C9").Value = 0
$ws1.Range("D9").Value = 1
$ws1.Range("E9").Value = 1
$ws1.Range("F9").Value = 0

$ws1.Range("A10").Value = 8
$ws1.Range("B10").Value = "D.Schultz"
$ws1.Range("C10").Value = 0
$ws1.Range("D10").Value = 1
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 0

# Copy style for new row 10 column A to match existing pattern
$ws1.Range("A9").Copy()
$ws1.Range("A10").PasteSpecial(-4122)

# --- Sheet2 (Receiving) ---
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "E.Elliott"
$ws2.Range("C2").Value = 40
$ws2.Range("D2").Value = 29
$ws2.Range("E2").Value = 2
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = 14
$ws2.Range("H2").Value = 10

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "T.Pollard"
$ws2.Range("C3").Value = 26
$ws2.Range("D3").Value = 22
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Value = 3
$ws2.Range("H3").Value = 3

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "C.Clement"
$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = 1
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 1
$ws2.Range("H4").Value = 1

$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "A.Cooper"
$ws2.Range("C5").Value = 52
$ws2.Range("D5").Value = 39
$ws2.Range("E5").Value = 21
$ws2.Range("F5").Value = 12
$ws2.Range("G5").Value = 12
$ws2.Range("H5").Value = 8

$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "C.Lamb"
$ws2.Range("C6").Value = 65
$ws2.Range("D6").Value = 45
$ws2.Range("E6").Value = 31
$ws2.Range("F6").Value = 16
$ws2.Range("G6").Value = 10
$ws2.Range("H6").Value = 5

$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "M.Gallup"
$ws2.Range("C7").Value = 38
$ws2.Range("D7").Value = 30
$ws2.Range("E7").Value = 10
$ws2.Range("F7").Value = 7
$ws2.Range("G7").Value = 4
$ws2.Range("H7").Value = 3

$ws2.Range("A8").Value = 6
$ws2.Range("B8").Value = "C.Wilson"
$ws2.Range("C8").Value = 22
$ws2.Range("D8").Value = 17
$ws2.Range("E8").Value = 9
$ws2.Range("F8").Value = 6
$ws2.Range("G8").Value = 5
$ws2.Range("H8").Value = 2

$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = "N.Brown"
$ws2.Range("C9").Value = 11
$ws2.Range("D9").Value = 8
$ws2.Range("E9").Value = 3
$ws2.Range("F9").Value = 2
$ws2.Range("G9").Value = 1
$ws2.Range("H9").Value = 1

$ws2.Range("A10").Value = 8
$ws2.Range("B10").Value = "M.Turner"
$ws2.Range("C10").Value = 8
$ws2.Range("D10").Value = 7
$ws2.Range("E10").Value = 2
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 4
$ws2.Range("H10").Value = 3

$ws2.Range("A11").Value = 9
$ws2.Range("B11").Value = "B.Jarwin"
$ws2.Range("C11").Value = 13
$ws2.Range("D11").Value = 10
$ws2.Range("E11").Value = 2
$ws2.Range("F11").Value = 0
$ws2.Range("G11").Value = 3
$ws2.Range("H11").Value = 2

$ws2.Range("A12").Value = 10
$ws2.Range("B12").Value = "D.Schultz"
$ws2.Range("C12").Value = 51
$ws2.Range("D12").Value = 39
$ws2.Range("E12").Value = 8
$ws2.Range("F12").Value = 5
$ws2.Range("G12").Value = 7
$ws2.Range("H12").Value = 3

$ws2.Range("A13").Value = 11
$ws2.Range("B13").Value = "S.McKeon"
$ws2.Range("C13").Value = 3
$ws2.Range("D13").Value = 1
$ws2.Range("E13").Value = 0
$ws2.Range("F13").Value = 0
$ws2.Range("G13").Value = 0
$ws2.Range("H13").Value = 0

# Copy style for new row 13 column A to match existing pattern
$ws2.Range("A12").Copy()
$ws2.Range("A13").PasteSpecial(-4122)

$excel.CutCopyMode = 0